# Driver Script file updated
# - Corrected the "Login" description text (typo fix: "execute" -> "execution")
# - Flipped Runmode from "Y" to "N" for every test case except Login_Verification
# - Left the active selection on B16 (where the user clicked after editing)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestSuite")

# Fix wording of the Login_Verification description
$ws.Range("B2").Value = "All type of login execution"

# Disable (Runmode = N) every test case below Login_Verification
$ws.Range("C3:C15").Value = "N"

# Leave the selection where the editor last clicked
$ws.Range("B16").Select() | Out-Null
